$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("汽車")
$ws1.Range("E2").Value = "98年03月24日"

$ws2 = $wb.Worksheets.Item("債務")
$ws2.Range("D2").Value = "合作金庫商業銀行臺南市北區曲門路"
$ws2.Range("F2").Value = "89年03月29日"

$wb.Save()
